$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Shift rows 66-68 down to 67-69 to make room for a new row 66 ---
# Work bottom-up so we don't clobber source data before it's copied.
$ws.Range("A68:O68").Copy($ws.Range("A69:O69"))
$ws.Range("A67:O67").Copy($ws.Range("A68:O68"))
$ws.Range("A66:O66").Copy($ws.Range("A67:O67"))

# --- 2) Populate the new row 66 with the FATF indicator ---
$ws.Range("A66").Value = "Z16_B04_P01_IB01"
$ws.Range("B66").Value = "Z16_B04_P01"
$ws.Range("C66").Value = "16.4"
$ws.Range("D66").Value = "Financial Action Task Force (FATF) rating effectiveness"
$ws.Range("E66").Value = "Financial Action Task Force (FATF) rating effectiveness"
$ws.Range("F66").Value = "XXX"
$ws.Range("G66").Value = "XXX"
$ws.Range("H66").Value = "XXX"
$ws.Range("I66").Value = "XXX"
$ws.Range("J66").Value = "XXX"
$ws.Range("K66").Value = "XXX"
$ws.Range("L66").Value = ""
$ws.Range("M66").Value = ""
$ws.Range("N66").Value = $false
$ws.Range("O66").ClearContents()

# --- 3) Update row 64 (Z16_B02_P01_Ib01) text fields ---
$ws.Range("D64").Value = "Von der Bundeswehr ausgebildetes Personal zur Stärkung der Kleinwaffenkontrolle und Munitionssicherheit"
$ws.Range("E64").Value = "XXXVon der Bundeswehr ausgebildetes Personal zur Stärkung der Kleinwaffenkontrolle und Munitionssicherheit"
$ws.Range("F64").Value = "XXX"
$ws.Range("G64").Value = "XXX"
$ws.Range("H64").Value = "XXX"
$ws.Range("I64").Value = "XXX"
$ws.Range("J64").Value = "XXX"
$ws.Range("K64").Value = "XXX"
$ws.Range("L64").Value = ""

# --- 4) Column width changes ---
$ws.Columns("C").ColumnWidth = 11.5703125
$ws.Columns("D").ColumnWidth = 26.07421875

# --- 5) Dimension will auto-expand to A1:O69 from the row writes above ---
